$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.954.98"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "1.831.87"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").Value = "'311.24"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").Value = "'0.4578"
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("D8").Value = "'0.3696"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.07186"
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("D10").Value = "'0.8781"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").Value = "'0.07814"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").Value = "'19.66"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "1.835.94"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("D14").Value = "'5.342"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").Value = "'6.412"
$ws.Range("E15").Value = "  -1.88%  "
$ws.Range("D16").Value = "'87.28"
$ws.Range("E16").Value = "  -4.95%  "
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").Value = "'0.000008725"
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "26.975.57"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("D21").Value = "'14.52"
$ws.Range("E21").Value = "  -1.98%  "
$ws.Range("D22").Value = "'5.014"
$ws.Range("E22").Value = "  -1.92%  "
$ws.Range("D23").Value = "2.057.70"
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").Value = "'10.44"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").Value = "'2.026"
$ws.Range("E25").Value = "  +7.27%  "
$ws.Range("D26").Value = "'151.46"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").Value = "'18.24"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("D28").Value = "'1.972"
$ws.Range("E28").Value = "  -4.69%  "
$ws.Range("D29").Value = "'114.17"
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("D30").Value = "'4.944"
$ws.Range("E30").Value = "  -3.37%  "
$ws.Range("D31").Value = "'0.08810"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").Value = "'3.033"
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("D33").Value = "'0.7553"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").Value = "'4.484"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "'1.137"
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("D36").Value = "'2.576"
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("D37").Value = "'1.090"
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("D38").Value = "'0.01938"
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("D39").Value = "'0.05156"
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("D40").Value = "'2.886"
$ws.Range("E40").Value = "  -3.48%  "
$ws.Range("D41").Value = "'6.954"
$ws.Range("E41").Value = "  -1.00%  "
$ws.Range("D42").Value = "'0.4989"
$ws.Range("E42").Value = "  -3.04%  "
$ws.Range("E43").Value = "  -2.24%  "
$ws.Range("D44").Value = "'8.342"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "'0.4696"
$ws.Range("E45").Value = "  -2.78%  "
$ws.Range("D46").Value = "'1.007"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("D47").Value = "'10.19"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").Value = "'102.41"
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").Value = "'1.616"
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("D50").Value = "'0.06113"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").Value = "'64.59"
$ws.Range("E51").Value = "  -1.42%  "
